# Update "Pais" sheet with refreshed country case data and re-sorted
# neighbours (country names swap position because their totals changed),
# plus the refreshed "datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Agosto de 2020 a las 23:36"

# --- Country name swaps (rows keep their position, but after refreshing
#     the totals the two neighbouring countries trade places) ----------
$ws.Range("A130").Value = "Ruanda"
$ws.Range("A131").Value = "Estonia"

$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# --- Refreshed numeric data --------------------------------------------
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 5352166
$ws.Range("C4").Value = 46209
$ws.Range("D4").Value = 2799525
$ws.Range("E4").Value = 2383774
$ws.Range("G4").Value = 1122
$ws.Range("H4").Value = 168867

# Row 5 - Brasil
$ws.Range("B5").Value = 3164785
$ws.Range("C5").Value = 52392
$ws.Range("E5").Value = 817460
$ws.Range("G5").Value = 1102
$ws.Range("H5").Value = 104201

# Row 13 - España
$ws.Range("H13").Value = 28579

# Row 22 - Alemania
$ws.Range("B22").Value = 220850
$ws.Range("C22").Value = 1320
$ws.Range("E22").Value = 11674
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = 9276

# Row 33 - Israel
$ws.Range("B33").Value = 88151
$ws.Range("C33").Value = 1558
$ws.Range("D33").Value = 62109
$ws.Range("E33").Value = 25403

# Row 53 - Barein
$ws.Range("B53").Value = 45264
$ws.Range("C53").Value = 460
$ws.Range("D53").Value = 41836
$ws.Range("E53").Value = 3262

# Row 92 - Gabon
$ws.Range("B92").Value = 8077
$ws.Range("C92").Value = 71
$ws.Range("D92").Value = 5920
$ws.Range("E92").Value = 2106

# Row 121 - Cabo Verde
$ws.Range("B121").Value = 3000
$ws.Range("C121").Value = 80
$ws.Range("D121").Value = 2172
$ws.Range("E121").Value = 795

# Row 124 - Mali
$ws.Range("B124").Value = 2582
$ws.Range("C124").Value = 5
$ws.Range("D124").Value = 1977
$ws.Range("E124").Value = 480

# Row 130 - now Ruanda (was Estonia's slot)
$ws.Range("B130").Value = 2189
$ws.Range("C130").Value = 18
$ws.Range("D130").Value = 1524
$ws.Range("E130").Value = 657
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 8

# Row 131 - now Estonia (was Ruanda's slot)
$ws.Range("B131").Value = 2174
$ws.Range("C131").Value = 7
$ws.Range("D131").Value = 1975
$ws.Range("E131").Value = 136
$ws.Range("H131").Value = 63

# Row 135 - Sierra Leona
$ws.Range("B135").Value = 1937
$ws.Range("C135").Value = 5
$ws.Range("D135").Value = 1483

# Row 136 - Yemen
$ws.Range("B136").Value = 1841
$ws.Range("C136").Value = 10
$ws.Range("D136").Value = 937
$ws.Range("E136").Value = 376
$ws.Range("G136").Value = 5
$ws.Range("H136").Value = 528

# Row 213 - now Islas Malvinas (was Montserrat's slot)
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214 - now Montserrat (was Islas Malvinas's slot)
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
